$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: B1 becomes "Correct Answer", and add 7 new
# "Wrong Answer N" headers across C1:I1 for the expanded game preview.
$ws.Range("B1").Value = "Correct Answer"
$ws.Range("C1").Value = "Wrong Answer 1"
$ws.Range("D1").Value = "Wrong Answer 2"
$ws.Range("E1").Value = "Wrong Answer 3"
$ws.Range("F1").Value = "Wrong Answer 4"
$ws.Range("G1").Value = "Wrong Answer 5"
$ws.Range("H1").Value = "Wrong Answer 6"
$ws.Range("I1").Value = "Wrong Answer 7"

# Move the active selection from H16 to F16.
$ws.Range("F16").Select()
